$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 231 (pushes the existing rows 231-264 down to 232-265,
# carrying their data and formatting with them).
$ws.Rows(231).Insert()

# Populate the newly inserted row 231 with the new weekly record.
$ws.Range("A231").Value = 10
$ws.Range("B231").Value = "Vega Modelo de Temuco"
$ws.Range("C231").Value = "La Araucanía"
$ws.Range("D231").Value = 44491
$ws.Range("E231").Value = 9
$ws.Range("F231").Value = 100114014
$ws.Range("G231").Value = "Betarraga"
$ws.Range("H231").Value = "Sin especificar"
$ws.Range("I231").Value = "Primera"
$ws.Range("J231").Value = 110
$ws.Range("K231").Value = 10000
$ws.Range("L231").Value = 10000
$ws.Range("M231").Value = 10000
$ws.Range("N231").Value = "$/docena de paquetes"
$ws.Range("O231").Value = "Provincia de Cautín"
$ws.Range("P231").Value = 833
$ws.Range("Q231").Value = 12
$ws.Range("R231").Value = "Hortaliza"
